$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameter")

# Update the informational text cells to reflect the new release
$ws.Range("A4").Value = "(You need a text parameter in ReportServer with key=cusNumber. You can then pass any customer name, e.g. 350)"
$ws.Range("A5").Value = "(Version: 1.0.1)"
$ws.Range("A6").Value = "(Last tested with: ReportServer 4.0.0-6053) "

# Move the active selection to A5
$ws.Range("A5").Select()
